$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 69 (hunk 0)
$ws.Range("H69").Value = 7299.095
$ws.Range("J69").Value = 7356.6343
$ws.Range("L69").Value = 22069.9029
$ws.Range("N69").Value = -23817.9029

# row 72 (hunk 1)
$ws.Range("H72").Value = 7299.095
$ws.Range("J72").Value = 7356.6343
$ws.Range("L72").Value = 66209.7087
$ws.Range("N72").Value = -74945.7087

# row 100 (hunk 2)
$ws.Range("H100").Value = 1425
$ws.Range("I100").Value = 850
$ws.Range("K100").Value = 850
$ws.Range("M100").Value = -309

$ws = $wb.Worksheets.Item("ARM")
# row 24 (hunk 3)
$ws.Range("H24").Value = 2520677.5
$ws.Range("J24").Value = 2520677.5
$ws.Range("L24").Value = 2520677.5
$ws.Range("N24").Value = -2521425.5

# row 46 (hunk 4)
$ws.Range("H46").Value = 3499
$ws.Range("I46").Value = 3499
$ws.Range("K46").Value = 3499
$ws.Range("M46").Value = -3180

# row 74 (hunk 5)
$ws.Range("H74").Value = 3705.0527
$ws.Range("I74").Value = 3429.4707
$ws.Range("K74").Value = 3429.4707
$ws.Range("M74").Value = -2555.4707

# row 77 (hunk 6)
$ws.Range("H77").Value = 3705.0527
$ws.Range("I77").Value = 3429.4707
$ws.Range("K77").Value = 17147.3535
$ws.Range("M77").Value = -12779.3535

# row 96 (hunk 7)
$ws.Range("H96").Value = 2886160.8
$ws.Range("J96").Value = 2886160.8
$ws.Range("L96").Value = 2886160.8
$ws.Range("N96").Value = -2891652.8

# row 100 (hunk 8)
$ws.Range("H100").Value = 2520677.5
$ws.Range("J100").Value = 2520677.5
$ws.Range("L100").Value = 2520677.5
$ws.Range("N100").Value = -2522841.5

# row 102 (hunk 9)
$ws.Range("H102").Value = 5955325
$ws.Range("I102").Value = 8334122
$ws.Range("K102").Value = 8334122
$ws.Range("M102").Value = -8332500

# row 103 (hunk 10)
$ws.Range("H103").Value = 47498.5
$ws.Range("J103").Value = 47498.5
$ws.Range("L103").Value = 47498.5
$ws.Range("N103").Value = -49842.5

$ws = $wb.Worksheets.Item("BSM")
# row 99 (hunk 11)
$ws.Range("H99").Value = 100001180
$ws.Range("I99").Value = 125001210
$ws.Range("J99").Value = 1107.5
$ws.Range("K99").Value = 125001210
$ws.Range("L99").Value = 1107.5
$ws.Range("M99").Value = -124999712
$ws.Range("N99").Value = -4103.5

# row 105 (hunk 12)
$ws.Range("H105").Value = 5683090
$ws.Range("I105").Value = 8265530.5
$ws.Range("K105").Value = 8265530.5
$ws.Range("M105").Value = -8263783.5

$ws = $wb.Worksheets.Item("CRP")
# row 141 (hunk 13)
$ws.Range("H141").Value = 75326
$ws.Range("J141").Value = 75326
$ws.Range("L141").Value = 75326
$ws.Range("N141").Value = -85686

$ws = $wb.Worksheets.Item("CUL")
# row 25 (hunk 14)
$ws.Range("H25").Value = 850
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 850
$ws.Range("K25").Value = 0
$ws.Range("L25").ClearContents()
$ws.Range("M25").Value = 2550
$ws.Range("N25").Value = -2888

# row 30 (hunk 15)
$ws.Range("H30").Value = 850
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 850
$ws.Range("K30").Value = 0
$ws.Range("L30").ClearContents()
$ws.Range("M30").Value = 2550
$ws.Range("N30").Value = -2754

# row 50 (hunk 16)
$ws.Range("H50").Value = 123.333336
$ws.Range("J50").Value = 123.333336
$ws.Range("L50").Value = 370.000008
$ws.Range("N50").Value = -1332.000008

# row 53 (hunk 17)
$ws.Range("H53").Value = 123.333336
$ws.Range("J53").Value = 123.333336
$ws.Range("L53").Value = 370.000008
$ws.Range("N53").Value = -1332.000008

# row 56 (hunk 18)
$ws.Range("H56").Value = 6665
$ws.Range("I56").Value = 6665
$ws.Range("K56").Value = 6665
$ws.Range("M56").Value = -6135

# row 69 (hunk 19)
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").ClearContents()
$ws.Range("N69").Value = 0

# row 72 (hunk 20)
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").ClearContents()
$ws.Range("N72").Value = 0

# row 104 (hunk 21)
$ws.Range("H104").Value = 8199.9

# row 112 (hunk 22)
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()

# row 131 (hunk 23)
$ws.Range("H131").Value = 3472.25
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 3472.25
$ws.Range("K131").Value = 0
$ws.Range("L131").ClearContents()
$ws.Range("M131").Value = 10416.75
$ws.Range("N131").Value = -20496.75

# row 137 (hunk 24)
$ws.Range("H137").Value = 3838.6667
$ws.Range("J137").Value = 3838.6667
$ws.Range("L137").Value = 11516.0001
$ws.Range("N137").Value = -21716.0001

$ws = $wb.Worksheets.Item("GSM")
# row 41 (hunk 25)
$ws.Range("H41").Value = 2000
$ws.Range("I41").Value = 2000
$ws.Range("K41").Value = 2000
$ws.Range("M41").Value = -1645

# row 97 (hunk 26)
$ws.Range("H97").Value = 1205.1
$ws.Range("I97").Value = 1410
$ws.Range("J97").Value = 1000.2
$ws.Range("K97").Value = 1410
$ws.Range("L97").Value = 1000.2
$ws.Range("M97").Value = -914
$ws.Range("N97").Value = -1992.2

# row 102 (hunk 27)
$ws.Range("H102").Value = 2172
$ws.Range("I102").Value = 2023.1578
$ws.Range("K102").Value = 2023.1578
$ws.Range("M102").Value = -401.1578

# row 122 (hunk 28)
$ws.Range("H122").Value = 1674.5
$ws.Range("I122").Value = 1081.1111
$ws.Range("K122").Value = 3243.3333
$ws.Range("M122").Value = -793.3333000000002

$ws = $wb.Worksheets.Item("LTW")
# row 93 (hunk 29)
$ws.Range("H93").Value = 4423
$ws.Range("I93").Value = 3843
$ws.Range("J93").Value = 5003
$ws.Range("K93").Value = 3843
$ws.Range("L93").Value = 5003
$ws.Range("M93").Value = -2595
$ws.Range("N93").Value = -7499

# row 99 (hunk 30)
$ws.Range("H99").Value = 50000
$ws.Range("J99").Value = 50000
$ws.Range("L99").Value = 50000
$ws.Range("N99").Value = -55990

# row 100 (hunk 31)
$ws.Range("H100").Value = 4569.154
$ws.Range("I100").Value = 4866.6665
$ws.Range("J100").Value = 3899.75
$ws.Range("K100").Value = 4866.6665
$ws.Range("L100").Value = 3899.75
$ws.Range("M100").Value = -4325.6665
$ws.Range("N100").Value = -4981.75

# row 122 (hunk 32)
$ws.Range("H122").Value = 2999.5
$ws.Range("I122").Value = 2999
$ws.Range("K122").Value = 8997
$ws.Range("M122").Value = -6547

$ws = $wb.Worksheets.Item("WVR")
# row 29 (hunk 33)
$ws.Range("H29").Value = 4755
$ws.Range("I29").Value = 4510
$ws.Range("K29").Value = 4510
$ws.Range("M29").Value = -4220

# row 80 (hunk 34)
$ws.Range("H80").Value = 38166.668
$ws.Range("J80").Value = 38166.668
$ws.Range("L80").Value = 38166.668
$ws.Range("N80").Value = -40162.668

# row 83 (hunk 35)
$ws.Range("H83").Value = 38166.668
$ws.Range("J83").Value = 38166.668
$ws.Range("L83").Value = 114500.004
$ws.Range("N83").Value = -124484.004

# row 96 (hunk 36)
$ws.Range("H96").Value = 10878
$ws.Range("I96").Value = 10053.6
$ws.Range("J96").Value = 15000
$ws.Range("K96").Value = 10053.6
$ws.Range("L96").Value = 15000
$ws.Range("M96").Value = -8680.6
$ws.Range("N96").Value = -17746
